$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (23, 24) need to inherit the bordered/bold style used by column A
# from the existing last data row (22) before we overwrite its contents.
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A23:A24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'2023-09-30"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = -3
$ws.Range("F2").Value = -3
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'2023-08-20"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 15
$ws.Range("E3").Value = -4
$ws.Range("F3").Value = -4
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'2023-08-19"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 9
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'2023-08-08"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 22
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 9

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'2023-07-15"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 30
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'2023-06-13"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 38
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 12
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 9

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'2023-05-06"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 24
$ws.Range("E8").Value = -4
$ws.Range("F8").Value = -3
$ws.Range("G8").Value = 8
$ws.Range("H8").Value = 2
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 5

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'2023-05-01"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 25
$ws.Range("E9").Value = -7
$ws.Range("F9").Value = -8
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 8

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'2023-04-30"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 19
$ws.Range("E10").Value = -1
$ws.Range("F10").Value = -1
$ws.Range("G10").Value = 8
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 8
$ws.Range("J10").Value = 8

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "'2023-04-26"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 31
$ws.Range("E11").Value = -2
$ws.Range("F11").Value = -1
$ws.Range("G11").Value = 9
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 5
$ws.Range("J11").Value = 6

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "'2023-04-23"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 19
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 7
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 9
$ws.Range("J12").Value = 9

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "'2023-04-20"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 25
$ws.Range("E13").Value = -4
$ws.Range("F13").Value = -5
$ws.Range("G13").Value = 4
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 8
$ws.Range("J13").Value = 9

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "'2023-04-18"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = -3
$ws.Range("F14").Value = -3
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 8
$ws.Range("J14").Value = 8

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "'2023-04-16"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 15
$ws.Range("E15").Value = 11
$ws.Range("F15").Value = 11
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 9

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "'2023-03-26"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -2
$ws.Range("F16").Value = -3
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 6

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "'2023-03-23"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 7

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "'2023-03-19"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 36
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 2
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 2

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "'2023-03-15"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = 6
$ws.Range("F19").Value = 6
$ws.Range("G19").Value = 9
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = 5

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "'2023-03-10"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 17
$ws.Range("E20").Value = 5
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 7

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "'2023-03-06"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = -2
$ws.Range("F21").Value = -2
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 5
$ws.Range("J21").Value = 5

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "'2023-03-03"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 11
$ws.Range("E22").Value = -3
$ws.Range("F22").Value = -3
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 9
$ws.Range("J22").Value = 9

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "'2023-02-28"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 11
$ws.Range("E23").Value = 5
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 6
$ws.Range("J23").Value = 6

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "'2023-02-25"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = -1
$ws.Range("F24").Value = -1
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 1
$ws.Range("I24").Value = 4
$ws.Range("J24").Value = 4

